$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     rows that are ready for handoff (rows 4-7) to the new report time.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-09-06 22:39:41"

# --- zh-cn sheet: mark rows 4-7 Priority as "ht" (was "low") and update
#     the Latest Handoff Datetime to the freshly generated report time.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-06 22:39:35"

# --- de-de sheet: same Priority update, and Latest Handoff Datetime
#     matches the Overview generate date for this handoff run.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-06 22:39:41"
